$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "angkush-2" -> "another one"
$ws.Range("B2").Value = "another one"

# D2: "8876690053" -> "1"
# The original value is stored as text (a phone number), so force the
# replacement to stay text as well instead of being auto-coerced to a
# number, then restore the default "Normal" style so no formatting
# residue is left behind.
$cellD2 = $ws.Range("D2")
$cellD2.NumberFormat = "@"
$cellD2.Value = "1"
$cellD2.Style = "Normal"

# E2: "hobby-2" -> "something"
$ws.Range("E2").Value = "something"
